$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.05105
$ws.Range("E2").Value = 0.0321
$ws.Range("G2").Value = 0.1790322580645161
$ws.Range("H2").Value = 0.1790322580645161
$ws.Range("I2").Value = 0.12
$ws.Range("J2").Value = 0.1082318241903503
$ws.Range("K2").Value = 48.7
$ws.Range("L2").Value = 0.1570967741935484
$ws.Range("M2").Value = 17.921
$ws.Range("N2").Value = 0.02902186234817814
$ws.Range("O2").Value = 0.3679876796714578
$ws.Range("P2").Value = 13.341
$ws.Range("Q2").Value = 0.02160485829959514
$ws.Range("R2").Value = 0.2739425051334702
$ws.Range("S2").Value = 4.58
$ws.Range("T2").Value = 0.2555660956419843
$ws.Range("U2").Value = 35
$ws.Range("V2").Value = 0.05668016194331984
$ws.Range("W2").Value = 0.2029011559614899
$ws.Range("X2").Value = 0.05765304454676771
$ws.Range("Y2").Value = 0.1452481114147222
$ws.Range("Z2").Value = 1.288124324773539
$ws.Range("AA2").Value = 0.135618051279294
$ws.Range("AB2").Value = 0.04698921527015996
$ws.Range("AC2").Value = 0.08862883600913406
$ws.Range("AD2").Value = 264.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 264.4
$ws.Range("AG2").Value = 229.4
$ws.Range("AH2").Value = 0.2998072343803153
$ws.Range("AI2").Value = 0.4842490842490843
$ws.Range("AJ2").Value = 0.2708702326130594
$ws.Range("AK2").Value = 0.4489236790606654
$ws.Range("AL2").Value = 2.41
$ws.Range("AM2").Value = 2.41
$ws.Range("AN2").Value = 6.448780487804878
$ws.Range("AO2").Value = 15.4356846473029
$ws.Range("AP2").Value = 5.595121951219513
$ws.Range("AQ2").Value = 15.4356846473029
$ws.Range("D3").Value = 0.0601
$ws.Range("E3").Value = -0.0103
$ws.Range("G3").Value = 0.2267857142857143
$ws.Range("H3").Value = 0.2267857142857143
$ws.Range("I3").Value = 0.1589285714285714
$ws.Range("J3").Value = 0.1351190476190476
$ws.Range("K3").Value = 22.7
$ws.Range("L3").Value = 0.1351190476190476
$ws.Range("M3").Value = 17.921
$ws.Range("N3").Value = 0.0624860529986053
$ws.Range("O3").Value = 0.7894713656387665
$ws.Range("P3").Value = 13.341
$ws.Range("Q3").Value = 0.04651673640167363
$ws.Range("R3").Value = 0.5877092511013216
$ws.Range("S3").Value = 4.58
$ws.Range("T3").Value = 0.2555660956419843
$ws.Range("U3").Value = 10.5
$ws.Range("V3").Value = 0.03661087866108786
$ws.Range("W3").Value = 0.1769290724863601
$ws.Range("X3").Value = 0.04615961871760139
$ws.Range("Y3").Value = 0.1307694537687587
$ws.Range("Z3").Value = 1.387741615727738
$ws.Range("AA3").Value = 0.1875103254584503
$ws.Range("AB3").Value = 0.04596228270772536
$ws.Range("AC3").Value = 0.141548042750725
$ws.Range("AD3").Value = 2.3
$ws.Range("AF3").Value = 2.3
$ws.Range("AG3").Value = -8.199999999999999
$ws.Range("AH3").Value = 0.007955724662746454
$ws.Range("AI3").Value = 0.01702442635085122
$ws.Range("AJ3").Value = -0.02943287867910983
$ws.Range("AK3").Value = -0.06581059390048154
$ws.Range("AN3").Value = 0.0812720848056537
$ws.Range("AP3").Value = -0.2897526501766784
$ws.Range("D4").Value = 0.042
$ws.Range("E4").Value = 0.0745
$ws.Range("G4").Value = 0.1225352112676056
$ws.Range("H4").Value = 0.1225352112676056
$ws.Range("I4").Value = 0.07394366197183098
$ws.Range("J4").Value = 0.07051833057166527
$ws.Range("K4").Value = 26
$ws.Range("L4").Value = 0.1830985915492958
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("T4").ClearContents()
$ws.Range("U4").Value = 24.5
$ws.Range("V4").Value = 0.07408527366192924
$ws.Range("W4").Value = 0.2288732394366197
$ws.Range("X4").Value = 0.06914647037593402
$ws.Range("Y4").Value = 0.1597267690606857
$ws.Range("Z4").Value = 1.187290969899666
$ws.Range("AA4").Value = 0.08372577710013772
$ws.Range("AB4").Value = 0.04801614783259458
$ws.Range("AC4").Value = 0.03570962926754314
$ws.Range("AD4").Value = 262.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 262.1
$ws.Range("AG4").Value = 237.6
$ws.Range("AH4").Value = 0.4421390013495277
$ws.Range("AI4").Value = 0.6378680944268679
$ws.Range("AJ4").Value = 0.4180890374802042
$ws.Range("AK4").Value = 0.6149068322981367
$ws.Range("AL4").Value = 2.41
$ws.Range("AM4").Value = 2.41
$ws.Range("AN4").Value = 20.63779527559056
$ws.Range("AO4").Value = 4.356846473029045
$ws.Range("AP4").Value = 18.70866141732284
$ws.Range("AQ4").Value = 4.356846473029045
